$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 52.049028
$ws.Range("H2").Value = 156.147084
$ws.Range("I2").Value = 0.7208330343078339
$ws.Range("J2").Value = 0.7208330343078339
$ws.Range("M2").Value = 1.672411
$ws.Range("N2").Value = 5.017233
$ws.Range("O2").Value = 0.6245395681653219
$ws.Range("P2").Value = 0.624539568165322
$ws.Range("Q2").Value = 87.04736696650801
$ws.Range("R2").Value = 783.426302698572
$ws.Range("S2").Value = 0.4501887519659132
$ws.Range("T2").Value = 0.4501887519659133
$ws.Range("G3").Value = 52.049028
$ws.Range("H3").Value = 156.147084
$ws.Range("I3").Value = 0.7208330343078339
$ws.Range("J3").Value = 0.7208330343078339
$ws.Range("O3").Value = 0.2513435317223857
$ws.Range("P3").Value = 0.2513435317223857
$ws.Range("Q3").Value = 35.03187589021601
$ws.Range("R3").Value = 315.286883011944
$ws.Range("S3").Value = 0.1811767206250946
$ws.Range("T3").Value = 0.1811767206250946
$ws.Range("G4").Value = 52.049028
$ws.Range("H4").Value = 156.147084
$ws.Range("I4").Value = 0.7208330343078339
$ws.Range("J4").Value = 0.7208330343078339
$ws.Range("N4").Value = 0.997092
$ws.Range("O4").Value = 0.1241169001122924
$ws.Range("P4").Value = 0.1241169001122924
$ws.Range("Q4").Value = 17.299223142192
$ws.Range("R4").Value = 155.693008279728
$ws.Range("S4").Value = 0.08946756171682604
$ws.Range("T4").Value = 0.08946756171682606
$ws.Range("I5").Value = 0.09317473454775864
$ws.Range("J5").Value = 0.09317473454775864
$ws.Range("M5").Value = 1.672411
$ws.Range("N5").Value = 5.017233
$ws.Range("O5").Value = 0.6245395681653219
$ws.Range("P5").Value = 0.624539568165322
$ws.Range("Q5").Value = 11.251725329117
$ws.Range("R5").Value = 101.265527962053
$ws.Range("S5").Value = 0.05819130847837568
$ws.Range("T5").Value = 0.05819130847837568
$ws.Range("I6").Value = 0.09317473454775864
$ws.Range("J6").Value = 0.09317473454775864
$ws.Range("O6").Value = 0.2513435317223857
$ws.Range("P6").Value = 0.2513435317223857
$ws.Range("S6").Value = 0.02341886684852944
$ws.Range("T6").Value = 0.02341886684852944
$ws.Range("I7").Value = 0.09317473454775864
$ws.Range("J7").Value = 0.09317473454775864
$ws.Range("N7").Value = 0.997092
$ws.Range("O7").Value = 0.1241169001122924
$ws.Range("P7").Value = 0.1241169001122924
$ws.Range("S7").Value = 0.01156455922085352
$ws.Range("T7").Value = 0.01156455922085352
$ws.Range("I8").Value = 0.1859922311444076
$ws.Range("J8").Value = 0.1859922311444076
$ws.Range("M8").Value = 1.672411
$ws.Range("N8").Value = 5.017233
$ws.Range("O8").Value = 0.6245395681653219
$ws.Range("P8").Value = 0.624539568165322
$ws.Range("Q8").Value = 22.460310816489
$ws.Range("R8").Value = 202.142797348401
$ws.Range("S8").Value = 0.116159507721033
$ws.Range("T8").Value = 0.1161595077210331
$ws.Range("I9").Value = 0.1859922311444076
$ws.Range("J9").Value = 0.1859922311444076
$ws.Range("O9").Value = 0.2513435317223857
$ws.Range("P9").Value = 0.2513435317223857
$ws.Range("Q9").Value = 9.039065148078002
$ws.Range("R9").Value = 81.35158633270201
$ws.Range("S9").Value = 0.04674794424876171
$ws.Range("T9").Value = 0.04674794424876171
$ws.Range("I10").Value = 0.1859922311444076
$ws.Range("J10").Value = 0.1859922311444076
$ws.Range("N10").Value = 0.997092
$ws.Range("O10").Value = 0.1241169001122924
$ws.Range("P10").Value = 0.1241169001122924
$ws.Range("S10").Value = 0.02308477917461283
$ws.Range("T10").Value = 0.02308477917461284